$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

# Update the "execute" flag values for the first two test cases (swap N/Y)
$ws.Range("D2").Value = "Y"
$ws.Range("D3").Value = "N"

# Update the selected cell on the sheet
$ws.Activate()
$ws.Range("E9").Select()
